$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (column F) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1137
$wsExhibit.Range("F4").Value = 2587
$wsExhibit.Range("F5").Value = 225

# Sheet "全部类型": update "想去人数" (column F) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1137
$wsAll.Range("F6").Value = 2587
$wsAll.Range("F8").Value = 225
